$wb = $excel.ActiveWorkbook

# --- 1. Text change: "Ready for handoff" -> "In Translation" ---
# This shared string is referenced from every sheet ("Overview", "zh-cn", "de-de");
# replace every occurrence so the underlying shared-string table only has
# one entry for it (matches the xlsx diff exactly).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        # NOTE: compare with the literal on the LEFT. PowerShell's -eq picks
        # its comparison semantics from the type of the LEFT operand, so
        # "$cell.Value2 -eq 'Ready for handoff'" would coerce the string to
        # boolean (and match) whenever the cell holds the boolean True.
        if ("Ready for handoff" -eq $cell.Value2) {
            $cell.Value = "In Translation"
        }
    }
}

# --- 2. Column width changes ---
# Overview sheet: columns E (5) and F (6) shrink from 17.2159881591797 to ~13.41
# zh-cn / de-de sheets: column C (3) shrinks the same way.
# (ColumnWidth set through COM is quantized to 1/6-character steps by the
# host, so 12.5 is the input that lands closest on the target stored width.)
$newColumnWidth = 12.5

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $newColumnWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newColumnWidth

$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Columns.Item(3).ColumnWidth = $newColumnWidth

$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Columns.Item(3).ColumnWidth = $newColumnWidth
